$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "24 décembre 2024"
$ws.Range("C3").Value = "Ouest Foire"
$ws.Range("D3").Value = "19h11"
$ws.Range("E3").Value = "Ben"
$ws.Range("G3").Value = "Yves"
$ws.Range("I3").Value = "SELUCY Taille Grande, SELUCY - VOSGIENNE Taille Petite"
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 8500
$ws.Range("L3").Value = 10000
